$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear previous data area below the header so no stale rows remain
$ws.Range("A2:C243").ClearContents()

$data = @(
  @("004431546", "GABRIELA", 397000),
  @("004480134", "JOSE", 62000),
  @("002697806", "CLAUDIA", 54916.67),
  @("001882235", "LAGO", 37343.33),
  @("005624730", "ISABEL", 30000.01),
  @("004377713", "DANIELI", 28672.04),
  @("005142661", "SABRINA", 11000.01),
  @("004313254", "GUSTAVO", 4292),
  @("004368468", "AHMAD", 3180.45),
  @("004213139", "LEONARDO", 2609.91),
  @("002823185", "SIMONE", 2165.46),
  @("008070407", "JULIANA", 1505.24),
  @("004398253", "EULER", 1398.44),
  @("001761119", "BLUEMETRIX", 1010.15),
  @("004482090", "CEZAR", 1003.01),
  @("004329030", "DANIELA", 940.23),
  @("004392159", "RODRIGO", 900.21),
  @("004458624", "PEDRO", 834.33),
  @("004211368", "ILTON", 826.8),
  @("004448303", "NASSIM", 774.03),
  @("005685353", "CARLOS", 767.05),
  @("005696595", "CLUBE", 752.05),
  @("004237325", "RICARDO", 606.19),
  @("004855960", "CLERIA", 556.35),
  @("004975924", "SERGIO", 532.85),
  @("005591536", "GUSTAVO", 520.17),
  @("004220849", "DULCE", 503.59),
  @("008002502", "JORGEANA", 500),
  @("000806386", "FERNANDA", 492.08),
  @("005002457", "ROSANGELA", 484.08),
  @("004488571", "CARLOS", 454.3),
  @("004432579", "ANA", 446.18),
  @("005186167", "ANDREA", 441.75),
  @("005044389", "CLAUDIA", 441.28),
  @("002687737", "JOSE", 370.35),
  @("004508516", "EDUARDO", 364.49),
  @("004355790", "MINEIA", 323.87),
  @("004207641", "MAGALI", 280.14),
  @("005040864", "ANDRE", 279.96),
  @("004374891", "RODRIGO", 273.15),
  @("004363260", "LARISSA", 257.51),
  @("004515341", "BRUNO", 235.12),
  @("003249855", "MARINA", 223.31),
  @("004472431", "LUIS", 197.96),
  @("003301389", "EDMUNDO", 191.02),
  @("004487016", "ROGERIO", 176.96),
  @("004508526", "CASSIO", 153.62),
  @("004224815", "GUILHERME", 133.6),
  @("004211911", "ZENILDA", 127.9),
  @("004381415", "JOAO", 127.55),
  @("008070544", "MARINA", 122.74),
  @("008090243", "GABRIEL", 122.74),
  @("004272426", "RODRIGO", 115.05),
  @("005135105", "BRENNER", 111.75),
  @("005646524", "EVANGELINA", 109.74),
  @("004207184", "CRISTINA", 100.15),
  @("008071998", "ISADORA", 100),
  @("005295509", "BHRUNA", 99.54),
  @("005701765", "F", 98.96),
  @("001368670", "THIAGO", 97.2),
  @("004239387", "LUIZ", 94.87),
  @("004431591", "MARIO", 94.24),
  @("004383268", "LAURA", 93.29),
  @("004384258", "PAULA", 93.29),
  @("004466350", "RAQUEL", 93.1),
  @("004536602", "TATIANY", 92.21),
  @("003115072", "VICTOR", 89.47),
  @("004462930", "WALTER", 87.61),
  @("004212132", "JOAO", 86.38),
  @("004809902", "PEDRO", 85.9),
  @("004517080", "TATIANA", 85.59),
  @("004261201", "ANA", 83.09),
  @("004424761", "PEDRO", 80),
  @("004384167", "DOUGLAS", 79.87),
  @("008069255", "ANGELA", 77.23),
  @("000827730", "LUCIANA", 76.01),
  @("005142611", "GUILHERME", 74.22),
  @("001719494", "LUIS", 73.48),
  @("004912314", "FABRICIO", 71.85),
  @("005880251", "LUIZ", 70.94),
  @("004563252", "FERNANDO", 70.58),
  @("004477812", "DIEGO", 70.23),
  @("004877741", "LUIZ", 70.02),
  @("005173958", "VENIA", 68.22),
  @("004452507", "DANIELA", 67.76),
  @("004212409", "RAFAEL", 67.39),
  @("004381194", "ALINNE", 67.06),
  @("008032257", "SARA", 67.03),
  @("004457389", "RAFAEL", 66.93),
  @("004749928", "NADY", 66.86),
  @("004242237", "MARIAH", 66.47),
  @("004027477", "GABRIELA", 64.77),
  @("004335251", "EDMUNDO", 62.82),
  @("005924958", "TIAGO", 62.24),
  @("003836362", "ISABELLA", 61.49),
  @("005092207", "BRUNO", 61.37),
  @("005068961", "JORGE", 61.09),
  @("005141215", "KARINA", 60.56),
  @("005890232", "TAYLA", 59.71),
  @("005558076", "ALEXANDRE", 59.47),
  @("004384131", "ANDRE", 58.71),
  @("004382374", "THEOMAR", 57.75),
  @("004974089", "CELIA", 57.68),
  @("004335144", "EDMUNDO", 57.28),
  @("004334062", "MERG", 57.01),
  @("004546050", "LUIS", 56.88),
  @("005079458", "JONAS", 56.38),
  @("004218542", "JOSE", 56.33),
  @("004813166", "VENIA", 55.69),
  @("004552021", "MARIA", 53.8),
  @("004504449", "KELMA", 53.39),
  @("004643880", "GABRIEL", 53.11),
  @("004400640", "FELIPE", 51.44),
  @("004426743", "GABRIELLE", 51.09),
  @("008032413", "VICTOR", 51.06),
  @("004332207", "IRACY", 48.95),
  @("005514036", "ANA", 48.75),
  @("005076418", "LEONARDO", 48.27),
  @("004693308", "LAURA", 48.25),
  @("005152037", "RODRIGO", 47.9),
  @("004260002", "ERICA", 46.9),
  @("004508159", "FELIPE", 46.28),
  @("001731007", "GUILHERME", 44.73),
  @("005981575", "GLAUCIANE", 44.17),
  @("004278212", "LEONARDO", 43.87),
  @("004381095", "THIAGO", 42.94),
  @("004332103", "JOSE", 42.34),
  @("004425965", "CAROLLINA", 41.94),
  @("004290978", "LARISSA", 40.34),
  @("004971448", "CLOVIS", 40.07),
  @("005245032", "ROSA", 39.91),
  @("004238164", "DANIELA", 39.09),
  @("004752519", "MARCUS", 37.58),
  @("004752615", "LUZIMAR", 37.11),
  @("004413537", "CLAUDIA", 36.55),
  @("005077648", "DUNAS", 36.46),
  @("004806286", "VERA", 35.77),
  @("005266369", "EG", 35.27),
  @("004240292", "MARCO", 34.71),
  @("005009992", "ALINE", 33.41),
  @("004435987", "MARCO", 32.34),
  @("004207374", "ANGELICA", 31.25),
  @("005055239", "NORMAN", 31.01),
  @("005070742", "JUSCELINO", 30.07),
  @("005927101", "SIMONE", 30),
  @("004332783", "IRON", 29.72),
  @("004377415", "ANGELA", 28.73),
  @("004230529", "LAIS", 28.18),
  @("004813134", "MONICA", 28.16),
  @("005305965", "SIDMAR", 28.12),
  @("005018038", "ELAINE", 27.47),
  @("004404724", "LEANDRO", 26.71),
  @("005616259", "MARIA", 26.65),
  @("004350197", "GISELA", 25.08),
  @("004472760", "SANDRA", 24.96),
  @("005715733", "ADRIANA", 24.5),
  @("004756968", "DANIELY", 24.09),
  @("004243043", "SUELI", 23.36),
  @("004398174", "DANIELE", 22.77),
  @("004371857", "NAZARETH", 21.52),
  @("004388077", "WLADMIR", 20.89),
  @("004214604", "MARIA", 20.75),
  @("004467884", "ANA", 20.69),
  @("005143579", "GABRIEL", 19.2),
  @("004204255", "AMADO", 18.77),
  @("004214592", "MERG", 18.64),
  @("004920447", "MARILIA", 18.19),
  @("008032597", "ALESSANDRO", 17.86),
  @("008037529", "MELISSA", 17.85),
  @("004994036", "BALTASAR", 17.62),
  @("003497496", "ELISANDRA", 17.15),
  @("002894447", "JOAO", 16.48),
  @("005293480", "WAGNER", 16.04),
  @("004547722", "MARCIA", 16.02),
  @("004340984", "RENATA", 15.62),
  @("004422594", "WANDIR", 14.67),
  @("004455356", "MARCELO", 14.59),
  @("004570632", "FABRICIO", 14.49),
  @("005878792", "JUNIO", 14.49),
  @("005268516", "LUIS", 14.27),
  @("004526450", "MSD", 14.13),
  @("004454365", "RAFAEL", 13.38),
  @("005685089", "CARNEIRO", 13.22),
  @("005009922", "ANA", 12.84),
  @("004264780", "MARCELO", 12.67),
  @("004374943", "LEONARDO", 12.08),
  @("004551472", "DIEGO", 11.88),
  @("004472076", "RUBENS", 11.8),
  @("005133039", "PAULO", 11.6),
  @("004335031", "EDMUNDO", 11.17),
  @("005374916", "MARCO", 10.89),
  @("004216298", "FLORDELIZ", 9.8),
  @("008013889", "CAROLINA", 9.72),
  @("002828327", "RENAN", 9.64),
  @("004751154", "CATARINE", 9.12),
  @("005324981", "JO", 7.69),
  @("004308815", "ZELI", 7.54),
  @("004505474", "RICARDO", 7.46),
  @("004530494", "ROSANGELA", 6.99),
  @("004752494", "SERGIO", 6.71),
  @("008007759", "CRISTINA", 6.46),
  @("004907688", "HEITOR", 6.12),
  @("004224405", "MILA", 5.88),
  @("004228456", "FLASH", 5.86),
  @("004448501", "JOAO", 5.55),
  @("004754056", "BRUNO", 5.18),
  @("004805269", "CLISIA", 5.18),
  @("005135281", "RAFAEL", 4.64),
  @("000834301", "MARCUS", 4.4),
  @("004539779", "RICARDO", 4.37),
  @("008012870", "ANA", 4.3),
  @("004165515", "MAURO", 4.2),
  @("008004995", "JOSE", 3.74),
  @("004328934", "VALERIA", 3.55),
  @("004352384", "BRASFORT", 3.13),
  @("004181486", "ANDREA", 3.11),
  @("005170415", "MONICA", 2.84),
  @("005140667", "MATEUS", 2.2),
  @("005022526", "ALEXANDRE", 1.7),
  @("000431814", "GUILHERME", 1.1),
  @("004360430", "VIOMAR", 1),
  @("004486497", "ELENA", 0.96),
  @("004115403", "HEBERT", 0.88),
  @("005660155", "CAROLINA", 0.85),
  @("004223502", "BRUNA", 0.78),
  @("004587511", "CARLOS", 0.69),
  @("004380948", "LUISA", 0.62),
  @("004473942", "DAIANNE", 0.62),
  @("003894173", "ANDREA", 0.48),
  @("004453302", "ISABELLA", 0.39),
  @("004638738", "GABRIEL", 0.29),
  @("004278033", "DAISY", 0.21),
  @("005662526", "AGUINALDO", 0.18),
  @("004432455", "LUCIANA", 0.17),
  @("002694089", "VITOR", 0.16),
  @("004357159", "JOAO", 0.15),
  @("004320840", "NATALIA", 0.14),
  @("001000288", "ISABELLA", 0.13),
  @("005530256", "CAROLINA", 0.1),
  @("004451996", "ADRIANO", 0.09),
  @("005047946", "GABRIEL", 0.09),
  @("004223226", "YESHUA", 0.04),
  @("004281300", "FRANKLIN", 0.04),
  @("005274028", "RAFAEL", 0.04),
  @("004329229", "GABRIEL", 0.03),
  @("004213373", "ALEXANDRE", 0.02),
  @("004339183", "JALISON", 0.02),
  @("004870976", "HFR", 0.02),
  @("000938440", "BASE", 0.01),
  @("002878817", "GUILHERME", 0.01),
  @("004400000", "VILMA", 0.01),
  @("004976625", "NORTON", 0.01),
  @("005324840", "PEDRO", -125.96)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = "'" + $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

$filterRow = $r + 1
$ws.Cells.Item($filterRow, 1).Value = "Filtros aplicados:`nDataFim é (Em branco)`nnr_saldo_disponivel não é 0`nPosição é Posição D-1`nCARTEIRA não está em branco`nDataFim é (Em branco)`nNR_CONTA não está em branco`nTIPO_LANCAMENTO não é ED, ET ou Liquidação Doador"